# Applies the lipidcane_spearman_1.xlsx edit:
#  - adds a new "Feedstock consumption [ton/yr]" metric column (column K),
#    shifting "Heat exchanger network error [%]" from column J to column K
#  - adds data in the newly-introduced column G ("Productivity [MMGGE/yr]")
#    for every data row (previously empty/missing)
#  - refreshes every Spearman correlation value in the C:K data block
#  - extends the header merge C1:J1 -> C1:K1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend header row formatting into the new column K -----------------
# Column K needs the same style (s="1") as the rest of the header cells in
# rows 1 and 2; copy formats from column J so we reuse the existing style
# instead of creating a brand-new one.
$ws.Range("J1:J2").Copy()
$ws.Range("K1:K2").PasteSpecial(-4122)  # xlPasteFormats

# --- Row 2: metric headers, with the new column inserted -----------------
$ws.Range("C2").Value = "MFPP [USD/ton]"
$ws.Range("D2").Value = "Biodiesel production [MMGal/yr]"
$ws.Range("E2").Value = "Ethanol production [MMGal/yr]"
$ws.Range("F2").Value = "Electricity production [MMWhr/yr]"
$ws.Range("G2").Value = "Natural gas consumption [MMcf/yr]"
$ws.Range("H2").Value = "Productivity [MMGGE/yr]"
$ws.Range("I2").Value = "TCI [10^6*USD]"
$ws.Range("J2").Value = "Feedstock consumption [ton/yr]"
$ws.Range("K2").Value = "Heat exchanger network error [%]"

# --- Data rows -------------------------------------------------------------
# Each row: B = parameter label, C..K = Spearman correlation coefficients.
# Column G (Productivity) is newly populated for every row.

$ws.Range("B4").Value  = "Lipid content [dry wt. %]"
$ws.Range("C4").Value  = -0.6692507312682818
$ws.Range("D4").Value  = 0.9455111377784444
$ws.Range("E4").Value  = -0.9746238655966399
$ws.Range("F4").Value  = -0.7387399684992125
$ws.Range("G4").Value  = 0.9717001008114807
$ws.Range("H4").Value  = -0.9712562814070353
$ws.Range("I4").Value  = 0.8680792019800497
$ws.Range("J4").Value  = -0.06124653116327909
$ws.Range("K4").Value  = 0.8636015900397511

$ws.Range("B5").Value  = "Lipid retention [%]"
$ws.Range("C5").Value  = -0.04949223730593266
$ws.Range("D5").Value  = 0.06121503037575941
$ws.Range("E5").Value  = -0.04262056551413786
$ws.Range("F5").Value  = -0.00648766219155479
$ws.Range("G5").Value  = 0.05310671232785925
$ws.Range("H5").Value  = -0.0487302182554564
$ws.Range("I5").Value  = 0.1414355358883972
$ws.Range("J5").Value  = 0.1361659041476037
$ws.Range("K5").Value  = -0.003388584714617866

$ws.Range("B6").Value  = "Bagasse lipid extraction efficiency [%]"
$ws.Range("C6").Value  = 0.09575189379734496
$ws.Range("D6").Value  = 0.2289492237305933
$ws.Range("E6").Value  = -0.01302932573314333
$ws.Range("F6").Value  = -0.1023730593264832
$ws.Range("G6").Value  = 0.0746671101052974
$ws.Range("H6").Value  = -0.02339908497712443
$ws.Range("I6").Value  = -0.005352133803345085
$ws.Range("J6").Value  = -0.06355658891472288
$ws.Range("K6").Value  = 0.257598439960999

$ws.Range("B7").Value  = "Capacity [ton/hr]"
$ws.Range("C7").Value  = 0.2084497112427811
$ws.Range("D7").Value  = 0.07732843321083029
$ws.Range("E7").Value  = 0.2422320558013951
$ws.Range("F7").Value  = 0.098775969399235
$ws.Range("G7").Value  = -0.1894033808159956
$ws.Range("H7").Value  = 0.2846921173029326
$ws.Range("I7").Value  = 0.2273846846171155
$ws.Range("J7").Value  = 0.7776224405610143
$ws.Range("K7").Value  = 0.05440336008400212

$ws.Range("B8").Value  = "Price [USD/gal]"
$ws.Range("C8").Value  = 0.404707117677942
$ws.Range("D8").Value  = 0.05014775369384235
$ws.Range("E8").Value  = -0.02640666016650417
$ws.Range("F8").Value  = -0.02898822470561765
$ws.Range("G8").Value  = 0.05632516182572641
$ws.Range("H8").Value  = -0.02670366759168979
$ws.Range("I8").Value  = 0.08715967899197483
$ws.Range("J8").Value  = 0.02759318982974574
$ws.Range("K8").Value  = 0.02388509712742819

$ws.Range("B9").Value  = "Price [USD/gal]"
$ws.Range("C9").Value  = 0.4829295732393311
$ws.Range("D9").Value  = -0.08734718367959199
$ws.Range("E9").Value  = 0.07976299407485186
$ws.Range("F9").Value  = 0.07667141678541965
$ws.Range("G9").Value  = -0.1064590552522034
$ws.Range("H9").Value  = 0.08326258156453914
$ws.Range("I9").Value  = -0.06393309832745819
$ws.Range("J9").Value  = -0.01898597464936624
$ws.Range("K9").Value  = -0.02834170854271357

$ws.Range("B10").Value = "Price [USD/cf]"
$ws.Range("C10").Value = 0.0346388659716493
$ws.Range("D10").Value = 0.01862296557413936
$ws.Range("E10").Value = -0.04279906997674943
$ws.Range("F10").Value = 0.01545938648466212
$ws.Range("G10").Value = 0.03014522538480182
$ws.Range("H10").Value = -0.03860196504912623
$ws.Range("I10").Value = 0.0005055126378159455
$ws.Range("J10").Value = -0.0102377559438986
$ws.Range("K10").Value = 0.02878271956798921

$ws.Range("B11").Value = "Electricity price [USD/kWh]"
$ws.Range("C11").Value = -0.02483162079051977
$ws.Range("D11").Value = 0.07964449111227782
$ws.Range("E11").Value = -0.03649891247281183
$ws.Range("F11").Value = 0.001512037800945024
$ws.Range("G11").Value = 0.06364587184619598
$ws.Range("H11").Value = -0.04004050101252531
$ws.Range("I11").Value = 0.1291982299557489
$ws.Range("J11").Value = 0.09863796594914874
$ws.Range("K11").Value = 0.04034200855021376

$ws.Range("B12").Value = "Operating days [day/yr]"
$ws.Range("C12").Value = 0.1043786094652366
$ws.Range("D12").Value = 0.1072541813545339
$ws.Range("E12").Value = 0.1292282307057677
$ws.Range("F12").Value = 0.05785344633615841
$ws.Range("G12").Value = -0.003881881669681422
$ws.Range("H12").Value = 0.05946898672466813
$ws.Range("I12").Value = -0.06612465311632792
$ws.Range("J12").Value = 0.5283537088427211
$ws.Range("K12").Value = -0.02064501612540313

$ws.Range("B13").Value = "IRR [%]"
$ws.Range("C13").Value = -0.2452381309532738
$ws.Range("D13").Value = -0.000475511887797195
$ws.Range("E13").Value = -0.008041701042526064
$ws.Range("F13").Value = -0.03633990849771244
$ws.Range("G13").Value = 0.01312861104015852
$ws.Range("H13").Value = -0.00301057526438161
$ws.Range("I13").Value = 0.05337283432085803
$ws.Range("J13").Value = 0.004314107852696318
$ws.Range("K13").Value = 0.006396159903997601

# --- Extend the merged header cell across the new column -----------------
# Merging re-derives per-cell borders for the merged block (splitting the
# box border into left/middle/right pieces), which would otherwise leave
# C1:K1 pointing at freshly-minted styles. Restore the original shared
# "header" style (same one already used by B1/B2/etc.) across the merged
# range afterwards so every cell keeps referencing the existing style index.
$ws.Range("C1:K1").Merge()
$ws.Range("B1").Copy()
$ws.Range("C1:K1").PasteSpecial(-4122)  # xlPasteFormats
